$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 1.9
$ws.Range("AJ2").Value = 16
$ws.Range("AO2").Value = 95
$ws.Range("F5").Value = 5.8
$ws.Range("I5").Value = 1.65
$ws.Range("L5").Value = 1.33
$ws.Range("N5").Value = 4.2
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 9.800000000000001
$ws.Range("AA5").Value = 16
$ws.Range("AJ5").Value = 230
$ws.Range("AL5").Value = 90
$ws.Range("G6").Value = 1.43
$ws.Range("I6").Value = 13.5
$ws.Range("K6").Value = 6
$ws.Range("R6").Value = 1.4
$ws.Range("S6").Value = 2.46
$ws.Range("V6").Value = 1.08
$ws.Range("F7").Value = 1.78
$ws.Range("G7").Value = 1.91
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 5.6
$ws.Range("L7").Value = 1.33
$ws.Range("N7").Value = 3.55
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 1.88
$ws.Range("Q7").Value = 1.94
$ws.Range("R7").Value = 1.34
$ws.Range("S7").Value = 3.05
$ws.Range("T7").Value = 1.84
$ws.Range("U7").Value = 1.96
$ws.Range("W7").Value = 2.08
$ws.Range("X7").Value = 17
$ws.Range("AC7").Value = 9.199999999999999
$ws.Range("AF7").Value = 11.5
$ws.Range("J9").Value = 3.2
$ws.Range("K9").Value = 3.7
$ws.Range("P9").Value = 1.76
$ws.Range("R9").Value = 1.29
$ws.Range("V9").Value = 1.37
$ws.Range("AN9").Value = 28
$ws.Range("G10").Value = 5.1
$ws.Range("P10").Value = 1.68
$ws.Range("H11").Value = 1.62
$ws.Range("L11").Value = 1.45
$ws.Range("R11").Value = 1.3
$ws.Range("AI11").Value = 44
$ws.Range("G13").Value = 4.7
$ws.Range("N13").Value = 3.05
$ws.Range("P13").Value = 1.69
$ws.Range("R13").Value = 1.26
$ws.Range("U13").Value = 1.88
$ws.Range("AG13").Value = 18.5
$ws.Range("N14").Value = 2.74
$ws.Range("P14").Value = 1.58
$ws.Range("Q14").Value = 2.54
$ws.Range("R14").Value = 1.21
$ws.Range("S14").Value = 5.2
$ws.Range("T14").Value = 2.06
$ws.Range("X14").Value = 8.800000000000001
$ws.Range("Y14").Value = 8.199999999999999
$ws.Range("F15").Value = 1.97
$ws.Range("H15").Value = 3.6
$ws.Range("I15").Value = 4.7
$ws.Range("J15").Value = 3.6
$ws.Range("L15").Value = 1.3
$ws.Range("Q15").Value = 1.76
$ws.Range("R15").Value = 1.41
$ws.Range("S15").Value = 2.94
$ws.Range("T15").Value = 1.68
$ws.Range("U15").Value = 2.2
$ws.Range("V15").Value = 1.31
$ws.Range("AA15").Value = 90
$ws.Range("AD15").Value = 19.5
$ws.Range("AF15").Value = 17
$ws.Range("AK15").Value = 26
$ws.Range("AO15").Value = 48
$ws.Range("Q19").Value = 1.82
$ws.Range("F21").Value = 1.47
$ws.Range("G21").Value = 1.61
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 9.6
$ws.Range("J21").Value = 3.6
$ws.Range("K21").Value = 5.2
$ws.Range("N21").Value = 2.88
$ws.Range("O21").Value = 1.26
$ws.Range("P21").Value = 1.92
$ws.Range("Q21").Value = 1.71
$ws.Range("R21").Value = 1.33
$ws.Range("S21").Value = 2.44
$ws.Range("V21").Value = 1.11
$ws.Range("W21").Value = 2.62

Write-Output "Applied 93 cell updates"
